$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Remove the now-duplicate chart defined names (_xlchart.v1.2 / v1.3)
$wb.Names.Item("_xlchart.v1.3").Delete()
$wb.Names.Item("_xlchart.v1.2").Delete()

# Add the two new header / formula columns (D & F) for mean/median increase
$ws.Range("D18").Font.Bold = $true
$ws.Range("F18").Font.Bold = $true
$ws.Range("D18").Value = "Mean increase"
$ws.Range("F18").Value = "Median increase"

$ws.Range("D19").Formula = "=((E3 / 95.321842) * 100) - 100"
$ws.Range("F19").Formula = "=((E10 / 95.22216) * 100) - 100"
$ws.Range("D19").Style = $ws.Range("A1").Style
$ws.Range("F19").Style = $ws.Range("A1").Style

# Update the active selection to match the author's saved view
$ws.Range("D21").Select()

$wb.Save()
